# Ability to discern between 2 conductor cables with attribute
#
# Sheet1 previously held 6 near-duplicate pull rows that only illustrated a
# single "STAR QUAD"/"7C#14" cable pulled between the same two stations.
# Replace that sample data with a single row describing a 2-conductor
# (2C#2) cable pulled between two new stations, and drop the now-irrelevant
# duplicate rows and the Distance figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 becomes the only data row: Cable Size -> 2C#2, From -> 100+00, To -> 200+00.
# Local/Express (column C) stays "EXPRESS".
$ws.Range("B2").Value = "2C#2"
$ws.Range("D2").Value = "100+00"
$ws.Range("E2").Value = "200+00"

# Distance is no longer tracked for this row.
$ws.Range("F2").ClearContents()

# Remove the obsolete duplicate rows (3-7), shifting everything below up.
$ws.Rows("3:7").Delete()

# Match the author's saved selection state.
$ws.Range("E3").Select()
